# Renaming the test ranges — refresh the r2 score matrix (C2:M18) with the
# newly computed values, and restore the Column A width that Excel recorded
# after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New r2 values for rows 2-18, columns C..M (11 columns per row).
$matrix = @(
    @(-0.55, -0.47, -0.62, -0.65, -0.36, -0.24, -0.76, -0.62, -0.34, -0.47, -0.51),
    @(-0.56, -0.46, -0.6,  -0.64, -0.34, -0.25, -0.76, -0.58, -0.34, -0.49, -0.49),
    @(-0.56, -0.46, -0.6,  -0.64, -0.35, -0.25, -0.76, -0.57, -0.34, -0.49, -0.49),
    @(-0.57, -0.47, -0.6,  -0.65, -0.35, -0.26, -0.77, -0.56, -0.34, -0.48, -0.49),
    @(-0.75, -0.52, -0.77, -0.47, -0.48, -0.63, -1.07, -0.94, -0.35, -0.61, -0.51),
    @(-2.6,  -0.52, -0.77, -0.59, -0.47, -1.25, -1.07, -0.94, -0.35, -0.61, -0.51),
    @(-1.9,  -0.52, -0.77, -0.89, -0.47, -0.31, -1.07, -0.94, -0.35, -0.61, -0.51),
    @(-0.74, -0.52, -0.77, -0.54, -0.47, -0.31, -1.07, -0.94, -0.35, -0.61, -0.51),
    @(-0.66, -0.28, -0.42, -0.22, -0.25, -0.42, -0.36, -0.28, -0.56, -0.38, -0.29),
    @(-0.6,  -0.29, -0.41, -0.17, -0.25, -0.43, -0.31, -0.25, -0.58, -0.38, -0.26),
    @(-0.6,  -0.3,  -0.41, -0.17, -0.26, -0.43, -0.31, -0.25, -0.58, -0.38, -0.26),
    @(-0.61, -0.3,  -0.42, -0.16, -0.26, -0.43, -0.31, -0.24, -0.58, -0.37, -0.25),
    @(-0.71, -0.35, -0.58, -0.55, -0.57, -0.33, -0.67, -0.62, -0.49, -0.45, -0.63),
    @(-1.41, -0.33, -0.58, -0.55, -0.56, -0.33, -0.67, -0.62, -0.49, -1.09, -0.63),
    @(-1.16, -0.33, -0.58, -0.55, -0.56, -0.33, -0.67, -0.62, -0.49, -0.5,  -0.63),
    @(-0.71, -0.32, -0.58, -0.55, -0.56, -0.33, -0.67, -0.62, -0.49, -0.49, -0.63),
    @(-0.56, -0.28, -0.41, -0.16, -0.25, -0.25, -0.31, -0.24, -0.34, -0.37, -0.25)
)

$startRow = 2
$startCol = 3   # column C

for ($i = 0; $i -lt $matrix.Length; $i++) {
    $rowValues = $matrix[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($r, $startCol + $j).Value = $rowValues[$j]
    }
}

# Column A got a bit wider after the edit.
$ws.Columns.Item(1).ColumnWidth = 21.75
